$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns G (ci.lower) and H (ci.upper) added to the results table
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

$ws.Range("G2").Value = -0.652956666624399
$ws.Range("H2").Value = -0.165726304037323
$ws.Range("G3").Value = -0.094656498555627
$ws.Range("H3").Value = -0.0164477837029356
$ws.Range("G4").Value = -0.0371802166385047
$ws.Range("H4").Value = 0.0282372982729989
$ws.Range("G5").Value = -0.646002269591595
$ws.Range("H5").Value = -0.112251200540699
$ws.Range("G6").Value = -0.328664743962842
$ws.Range("H6").Value = 0.249611359108821
$ws.Range("G7").Value = -0.729103838385689
$ws.Range("H7").Value = -0.126691166626034
$ws.Range("G8").Value = -0.286385394313043
$ws.Range("H8").Value = 0.217501416919475
$ws.Range("G9").Value = -0.0627142793191626
$ws.Range("H9").Value = 0.00269067900712837
$ws.Range("G10").Value = -0.150220010959749
$ws.Range("H10").Value = -0.0893147909789818
$ws.Range("G11").Value = -0.0814808002803514
$ws.Range("H11").Value = -0.0206805636127053
$ws.Range("G15").Value = -0.0180893926313762
$ws.Range("H15").Value = 0.0602271562123987
$ws.Range("G16").Value = -0.139335870255663
$ws.Range("H16").Value = 0.463907406671171
$ws.Range("G17").Value = -0.230428573379141
$ws.Range("H17").Value = 0.490975231164257
$ws.Range("G18").Value = -0.0743604332770094
$ws.Range("H18").Value = 0.0564745965459978
$ws.Range("G19").Value = -0.572770788626086
$ws.Range("H19").Value = 0.435002833838949
$ws.Range("G20").Value = -0.657329487925683
$ws.Range("H20").Value = 0.499222718217642
$ws.Range("G21").Value = 0.0206805636127053
$ws.Range("H21").Value = 0.0814808002803515
$ws.Range("G22").Value = -0.0053813580142567
$ws.Range("H22").Value = 0.125428558638325
$ws.Range("G23").Value = 0.159294697565236
$ws.Range("H23").Value = 0.627616330052918
$ws.Range("G24").Value = -0.0414506013194265
$ws.Range("H24").Value = 0.966129583724718
$ws.Range("G25").Value = 0.107566407734587
$ws.Range("H25").Value = 0.571633677543686
$ws.Range("G26").Value = -0.0874029028098363
$ws.Range("H26").Value = 0.92470975779615

# Rows 12-14 have no CI figures (matches blank C:F cells there already);
# touch the border property so Excel still materializes empty <c> cells
# for G/H on those rows without introducing any new cell style.
$ws.Range("G12:H14").Borders.LineStyle = -4142
